$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Drop the existing hyperlinks before we start moving columns around (their
#    anchors won't track the shift automatically).
$ws.Range("A1").Hyperlinks.Delete()

# 2. Remove the "Mac_Address" column (column C). This shifts D..H left to C..G,
#    putting IP_Address / Username / Password / Game_Executable / Processes
#    into their new positions.
$ws.Range("C1").EntireColumn.Delete()

# 3. Rename the shifted headers.
$ws.Range("C1").Value = "IP_Address(optional)"
$ws.Range("G1").Value = "Inserted_Processes(seperate by comma)"

# 4. Update data values.
$ws.Range("G2").Value = "mikesunique.exe,miketwo.exe"

$ws.Range("D3").Value = "Ed Greenlee"
$ws.Range("E3").Value = 2447

# 5. Re-attach the hyperlinks on the cells that now sit in column D, while the
#    cell text still matches the mail address, so the cached "display" text
#    comes out right; re-apply the built-in Hyperlink style (Add() style the
#    cell with a duplicate format otherwise); then overwrite row 4's
#    Username/Password with the dummy placeholder text used elsewhere in that
#    row.
$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:mmcquad.17@gmail.com")
$ws.Range("D2").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("D4"), "mailto:mmcquad.17@gmail.com", $null, $null, "mmcquad.17@gmail.com")
$ws.Range("D4").Style = "Hyperlink"

$ws.Range("D4").Value = "dummyMike"
$ws.Range("E4").Value = "dummyMike"
